# Enhance order management: update SQL queries for discounts, refactor order
# item retrieval, and improve UI for displaying product details.
#
# Concretely, this adds a new "Order_Info" column (G) to the Orders sheet,
# back-fills it for the existing rows, fixes up a total amount, and appends
# 12 new order rows that were captured after the refactor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header: G1 = "Order_Info" ------------------------------
# Copy the formatting of the last existing header cell (F1 - bold, bordered)
# onto the new header cell before setting its text.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Order_Info"

# --- Back-fill Order_Info for existing rows 2-31 ------------------------
for ($i = 1; $i -le 30; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 7).Value = "Random Order Info $i"
}

# --- Correction to an existing total amount -----------------------------
$ws.Range("E31").Value = 60.06

# --- Append newly captured orders (rows 32-43) --------------------------
$newRows = @(
  @(31, 2, 1, "2025-03-03 17:38:18", 281.4798305084746, 0),
  @(32, 2, 1, "2025-03-03 17:44:09", 75.075, 0),
  @(33, 2, 1, "2025-03-03 17:45:14", 0, 0),
  @(34, 2, 1, "2025-03-03 17:45:56", 55422.42177503629, 0),
  @(35, 2, 1, "2025-03-03 17:49:17", 99.46107784431139, 0),
  @(36, 2, 1, "2025-03-03 17:50:02", 11186.68115876002, 0),
  @(37, 2, 1, "2025-03-03 17:53:52", 2974.681818181818, 0),
  @(38, 2, 1, "2025-03-03 17:55:12", 1793.72197309417, 0),
  @(39, 2, 1, "2025-03-03 18:01:42", 269.298108982036, 0),
  @(40, 2, 1, "2025-03-03 18:05:52", 13071.4896021017, 0),
  @(41, 2, 1, "2025-03-03 18:12:41", 1883.81197309417, 0),
  @(42, 8, 2, "2025-03-03 18:19:56", 0, 0)
)

$startRow = 32
foreach ($r in $newRows) {
  $ws.Cells.Item($startRow, 1).Value = $r[0]
  $ws.Cells.Item($startRow, 2).Value = $r[1]
  $ws.Cells.Item($startRow, 3).Value = $r[2]
  $ws.Cells.Item($startRow, 4).Value = $r[3]
  $ws.Cells.Item($startRow, 5).Value = $r[4]
  $ws.Cells.Item($startRow, 6).Value = $r[5]
  $startRow++
}
